$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab from "ShearF-HW30.xpc" to "ShearF"
$ws.Name = "ShearF"

# New row 16 mirrors the formatting of row 15's first column (bold,
# centered, bordered style) - copy it over before writing values.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16").PasteSpecial(-4122) | Out-Null

# Populate the new data row (Gaussian-quadrature averaged intensities)
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9528535474453399
$ws.Range("D16").Value = 1.1851820878893
$ws.Range("E16").Value = 0.9489456796225956
$ws.Range("F16").Value = 0.9528535474453399
$ws.Range("G16").Value = 1.102144700832245
$ws.Range("H16").Value = 0.8620159336131725
$ws.Range("I16").Value = 0.9450172605927911
$ws.Range("J16").Value = 1.1851820878893
$ws.Range("K16").Value = 1.067063883755948
$ws.Range("L16").Value = 1.009958715600644
$ws.Range("M16").Value = 0.9993598683325741
